# Reproduce the author's edits to Sheet1 of the workbook:
#   - D5: 2 -> 3
#   - F5: 2 -> 3
#   - H5: 36 -> 46
#   - active selection ends on C5 (was D2)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").Value = 3
$ws.Range("F5").Value = 3
$ws.Range("H5").Value = 46

# Leave the selection where the author left it when they saved the file.
$ws.Range("C5").Select()
